$wb = $excel.ActiveWorkbook
Write-Host "Workbooks count: $($excel.Workbooks.Count)"
try {
  Write-Host "Name: $($excel.Workbooks.Item(1).Name)"
} catch { Write-Host "err1: $_" }
try {
  $links = $wb.LinkSources(1)
  Write-Host "links: $links"
} catch { Write-Host "err2: $_" }
